$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C5").Value = "Verificou lista de carros em produção"
